$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells G1, H1 -- reuse the exact formatting of the existing
# header cell F1 (bold, bordered, centered) via copy/paste-special formats.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("H1").PasteSpecial(-4122)

$ws.Range("G1").Value = "dim"
$ws.Range("H1").Value = "time_elapsed"

# Change C2 from a numeric timestamp to a formatted date string
$ws.Range("C2").Value = "Wednesday, January 1, 2020 00:00:00"

# Add new data cells G2, H2 (plain numeric, unstyled like F2)
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 1
